$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B5").Value = -0.31    # Total P&L %
$wsSummary.Range("B6").Value = 7        # Total Trades
$wsSummary.Range("B9").Value = 28.57    # Win Rate %

# --- Strategy Status sheet ---
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("D4").Value = 7         # MarketMaking Trades
$wsStatus.Range("G4").Value = 28.57     # MarketMaking Win Rate %

# --- New trade row (Trade #7) to append to "All Trades" and "MarketMaking" sheets ---
# Note: the date column (B) is a text value that looks like a date
# ("2026-02-17"); a leading apostrophe forces it to be stored as text
# instead of being auto-converted into a date serial number.
$newRow = @(7, "'2026-02-17", "15:13:57", "MarketMaking", "UP", 0.94, 0.94, "CLOSED", 0, 0, 99.89, 0, 0, 0.6, "Normal spread capture: 19600 bps", "early_exit", 0.14)

$wsAllTrades = $wb.Worksheets.Item("All Trades")
for ($i = 0; $i -lt $newRow.Length; $i++) {
    $wsAllTrades.Cells.Item(8, $i + 1).Value = $newRow[$i]
}

$wsMarketMaking = $wb.Worksheets.Item("MarketMaking")
for ($i = 0; $i -lt $newRow.Length; $i++) {
    $wsMarketMaking.Cells.Item(8, $i + 1).Value = $newRow[$i]
}
